$wb = $excel.ActiveWorkbook

# --- 1. Sheet1 (the monthly data sheet): insert a new header row above the
#        existing data and label the Target/Total columns. ---
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(1).Insert()
$ws.Range("E1").Value = "Target"
$ws.Range("F1").Value = "Total"

# Move the active selection on the data sheet to match the edited workbook.
$ws.Range("F1").Select()

# --- 2. Chart1: update the two line-chart series so their names and source
#        ranges follow the row shift caused by the inserted row. ---
$chartSheet = $wb.Worksheets.Item("Chart1")
$co = $chartSheet.ChartObjects().Item(1)
$chart = $co.Chart

$ser1 = $chart.SeriesCollection().Item(1)
$ser1.Name = "Target"
$ser1.Formula = '=SERIES("Target",Sheet1!$A$2:$A$13,Sheet1!$E$2:$E$13,1)'

$ser2 = $chart.SeriesCollection().Item(2)
$ser2.Name = "Actual"
$ser2.Formula = '=SERIES("Actual",Sheet1!$A$2:$A$13,Sheet1!$F$2:$F$13,2)'

# --- 3. Chart1 sheet view: rezoom and move the selection. ---
$chartSheet.Select()
$excel.ActiveWindow.Zoom = 90
$chartSheet.Range("Q33").Select()
